# Alarm Normal load method changes
# Add two new columns (I: AlarmLoadingDetail, J: StandbyLoadingDetail) to the
# "Add Panels" sheet, with corresponding per-row detail values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Add Panels")

# --- Header row (row 7): copy formatting from existing header cell H7 ---
$ws.Range("H7").Copy()
$ws.Range("I7:J7").PasteSpecial(-4122)
$ws.Range("I7").Value = "AlarmLoadingDetail"
$ws.Range("J7").Value = "StandbyLoadingDetail"

# --- Data rows (8-10): copy formatting from existing data cell H8 ---
$ws.Range("H8").Copy()
$ws.Range("I8:J10").PasteSpecial(-4122)

$ws.Range("I8").Value = "Battery Alarm (A)"
$ws.Range("J8").Value = "Battery Standby (A)"

$ws.Range("I9").Value = "Battery Alarm (A)"
$ws.Range("J9").Value = "Battery Standby (A)"

$ws.Range("I10").Value = "Battery Alarm (A)"
$ws.Range("J10").Value = "Battery Standby (A)"

# Restore the active selection to the newly added range, matching the
# author's last selection when saving the workbook.
$ws.Activate()
$ws.Range("I7:J10").Select()
